# Apply updated crypto price/volume data to Sheet1 (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, ForceText (avoid Excel turning numeric-looking
# strings like "245.04" into real numbers, which would also drop trailing zeros).
$updates = @(
    @{ Cell = 'D2'; Value = '41.810.38'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -4.83%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '2.216.38'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -5.94%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  +0.08%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '245.04'; ForceText = $true },
    @{ Cell = 'D6'; Value = '0.626'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  -6.58%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '68.99'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -7.01%  '; ForceText = $false },
    @{ Cell = 'E8'; Value = '  +0.15%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.546'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -8.96%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  -5.52%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '57.99'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -3.96%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '35.51'; ForceText = $true },
    @{ Cell = 'E12'; Value = '  +6.15%  '; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -2.96%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '6.69'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -8.08%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '2.546.77'; ForceText = $false },
    @{ Cell = 'E15'; Value = '  -5.78%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '14.78'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -8.90%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '0.839'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -7.67%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '2.216.54'; ForceText = $false },
    @{ Cell = 'E18'; Value = '  -5.59%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '41.770.08'; ForceText = $false },
    @{ Cell = 'E19'; Value = '  -4.74%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '0.0₃0954'; ForceText = $false },
    @{ Cell = 'E20'; Value = '  -7.18%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '72.52'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -7.30%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '6.03'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -9.24%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '234.70'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -7.20%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '2.04'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +9.64%  '; ForceText = $false },
    @{ Cell = 'E25'; Value = '  -0.17%  '; ForceText = $false },
    @{ Cell = 'B26'; Value = 'PancakeSwap'; ForceText = $false },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.46'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = 'B27'; Value = 'WEMIXToken'; ForceText = $false },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText = $false },
    @{ Cell = 'D27'; Value = '3.61'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -5.07%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '2.23'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -3.28%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '9.86'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  -5.78%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '170.59'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -3.08%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '20.37'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -8.64%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '0.120'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -6.09%  '; ForceText = $false },
    @{ Cell = 'E33'; Value = '  -7.59%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  -4.32%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '5.15'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -4.51%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  -8.12%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '3.86'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +0.98%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '22.66'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +16.80%  '; ForceText = $false },
    @{ Cell = 'E39'; Value = '  -5.27%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '0.0274'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +0.31%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '5.83'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -9.42%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '65.74'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +0.37%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '4.94'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -11.28%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '8.92'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -2.70%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '0.100'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  -5.18%  '; ForceText = $false },
    @{ Cell = 'B46'; Value = 'BinanceUSD'; ForceText = $false },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; ForceText = $false },
    @{ Cell = 'D46'; Value = '1.00'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +0.16%  '; ForceText = $false },
    @{ Cell = 'B47'; Value = 'Algorand'; ForceText = $false },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false },
    @{ Cell = 'D47'; Value = '0.188'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -6.50%  '; ForceText = $false },
    @{ Cell = 'E48'; Value = '  +7.66%  '; ForceText = $false },
    @{ Cell = 'B49'; Value = 'Celestia'; ForceText = $false },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; ForceText = $false },
    @{ Cell = 'D49'; Value = '10.11'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  +6.75%  '; ForceText = $false },
    @{ Cell = 'B50'; Value = 'TrustWalletToken'; ForceText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false },
    @{ Cell = 'D50'; Value = '1.17'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -4.49%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '1.10'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -4.94%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to treat the value as literal text
        # instead of coercing it into a Double (which would lose formatting
        # like trailing zeros, e.g. '1.00' -> 1).
        $range.Value = "'" + $u.Value
        # Setting Style back to Normal clears the quote-prefix flag that Excel
        # attaches to the cell style when a leading apostrophe is used, so the
        # cell ends up with the same (default) style as before the edit.
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
